# Update the cryptocurrency price/volume snapshot cells (columns D and E)
# with the latest scraped values from the GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.583.05"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.847.83"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'312.39"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4242"
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("D8").Value = "'0.3637"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D9").Value = "'44.41"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.07285"
$ws.Range("D11").Value = "'0.8725"
$ws.Range("E11").Value = "  -3.47%  "
$ws.Range("D12").Value = "'20.52"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "1.832.20"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").Value = "'5.322"
$ws.Range("E14").Value = "  -0.34%  "
$ws.Range("D15").Value = "'6.518"
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").Value = "'0.06898"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'79.78"
$ws.Range("E18").Value = "  +3.06%  "
$ws.Range("D19").Value = "'0.000008941"
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("D21").Value = "'15.34"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").Value = "27.617.12"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'4.982"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "'10.35"
$ws.Range("E24").Value = "  -2.92%  "
$ws.Range("D25").Value = "2.062.99"
$ws.Range("E25").Value = "  +2.26%  "
$ws.Range("D26").Value = "'1.980"
$ws.Range("E26").Value = "  -1.96%  "
$ws.Range("D27").Value = "'154.24"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "'18.83"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("D29").Value = "'121.96"
$ws.Range("E29").Value = "  +10.23%  "
$ws.Range("D30").Value = "'5.258"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "'1.876"
$ws.Range("E31").Value = "  +13.33%  "
$ws.Range("D32").Value = "'0.08875"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").Value = "'0.7653"
$ws.Range("E33").Value = "  -1.42%  "
$ws.Range("D34").Value = "'2.977"
$ws.Range("E34").Value = "  +2.02%  "
$ws.Range("D35").Value = "'4.535"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").Value = "'1.103"
$ws.Range("E36").Value = "  +3.49%  "
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("D38").Value = "'0.05366"
$ws.Range("E38").Value = "  +0.50%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("D40").Value = "'2.821"
$ws.Range("E40").Value = "  -4.45%  "
$ws.Range("D41").Value = "'6.887"
$ws.Range("E41").Value = "  +0.95%  "
$ws.Range("D42").Value = "'0.5069"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("E43").Value = "  +1.35%  "
$ws.Range("D44").Value = "'8.337"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "'0.06533"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("D46").Value = "'10.38"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").Value = "'0.4676"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "'104.66"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").Value = "'64.25"
$ws.Range("E51").Value = "  -0.19%  "
